# resource_log.xlsx edit:
# The Power Query-backed "resource_log" table is refreshed and now returns
# three rows (clear / train / process) instead of one (process). Row 2 is
# overwritten in place with the "clear" stage, and two more rows (train,
# process) are appended as rows 3-4. The calculated "duration_days" column
# (M) and the table/defined-name ranges grow to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("resource_log")

# --- Row 2: "clear" (overwrites the former single "process" row) ---
$ws.Range("A2").Value = "clear"
$ws.Range("B2").Value = 45796.46444909722
$ws.Range("C2").Value = 45796.464558368054
$ws.Range("D2").Value = 9.4416860000000007
$ws.Range("E2").Value = 4.9000000000000004
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = 56.1
$ws.Range("H2").Value = 55.9
$ws.Range("I2").NumberFormat = "General"
$ws.Range("J2").NumberFormat = "General"
$ws.Range("K2").NumberFormat = "General"
$ws.Range("L2").NumberFormat = "General"

# --- Row 3: "train" ---
$ws.Range("A3").Value = "train"
$ws.Range("B3").Value = 45796.464672835646
$ws.Range("C3").Value = 45796.491494606482
$ws.Range("D3").Value = 2317.4009489999999
$ws.Range("E3").Value = 21
$ws.Range("F3").Value = 3.3
$ws.Range("G3").Value = 58.1
$ws.Range("H3").Value = 41.7
$ws.Range("I3").NumberFormat = "General"
$ws.Range("J3").NumberFormat = "General"
$ws.Range("K3").NumberFormat = "General"
$ws.Range("L3").NumberFormat = "General"
$ws.Range("M3").Formula = '=CONVERT(resource_log[[#This Row],[duration_s]],"s","d")'

# --- Row 4: "process" ---
$ws.Range("A4").Value = "process"
$ws.Range("B4").Value = 45796.491508055558
$ws.Range("C4").Value = 45796.493430659721
$ws.Range("D4").Value = 166.11286799999999
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 6.5
$ws.Range("G4").Value = 41.9
$ws.Range("H4").Value = 51.1
$ws.Range("I4").NumberFormat = "General"
$ws.Range("J4").NumberFormat = "General"
$ws.Range("K4").NumberFormat = "General"
$ws.Range("L4").NumberFormat = "General"
$ws.Range("M4").Formula = '=CONVERT(resource_log[[#This Row],[duration_s]],"s","d")'

# Apply matching date/time formats to the new start/end timestamp cells so
# they render the same way as the existing row (numFmtId 22).
$ws.Range("B3:C4").NumberFormat = $ws.Range("B2").NumberFormat

# Grow the table (and its AutoFilter) to cover the new rows.
$lo = $ws.ListObjects.Item("resource_log")
$lo.Resize($ws.Range("A1:M4"))

# Grow the hidden Power Query "ExternalData_1" defined name to match the
# refreshed source range (columns A:L, 4 rows incl. header).
$wb.Names.Item("resource_log!ExternalData_1").RefersTo = "=resource_log!`$A`$1:`$L`$4"

# Recalculate so the M column formulas carry fresh cached values.
$excel.CalculateFull()

# Match the saved selection left behind by the edit.
$ws.Range("D2").Select()
